# Update odds values in the "Jogos da Semana" worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 changes
$ws.Range("K2").Value = 1.8
$ws.Range("Q2").Value = 3.5
$ws.Range("R2").Value = 1.3
$ws.Range("AG2").Value = 8.5
$ws.Range("AN2").Value = 3.25
$ws.Range("AR2").Value = 101
$ws.Range("AT2").Value = 13

# Row 3 changes
$ws.Range("G3").Value = 3.05
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 2.27
$ws.Range("J3").Value = 3.5
$ws.Range("L3").Value = 2.82
$ws.Range("P3").Value = 2.95
$ws.Range("Q3").Value = 1.93
$ws.Range("R3").Value = 1.78
$ws.Range("W3").Value = 9.5
$ws.Range("X3").Value = 16.5
$ws.Range("Y3").Value = 10.75
$ws.Range("AA3").Value = 27
$ws.Range("AI3").Value = 8.75
$ws.Range("AJ3").Value = 24
$ws.Range("AK3").Value = 18.5
$ws.Range("AL3").Value = 27
$ws.Range("AO3").Value = 16.5
$ws.Range("AQ3").Value = 75
$ws.Range("AR3").Value = 100
$ws.Range("AV3").Value = 4.25
$ws.Range("AW3").Value = 11.75
$ws.Range("AY3").Value = 45
$ws.Range("AZ3").Value = 70

$wb.Save()
